$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("M2").Value = 3.456265333333333
$ws.Range("N2").Value = 10.368796
$ws.Range("O2").Value = 0.009841535807677501
$ws.Range("P2").Value = 0.0098415358076775
$ws.Range("Q2").Value = 0.8128076142631111
$ws.Range("R2").Value = 7.315268528368
$ws.Range("S2").Value = 0.009586633444024011
$ws.Range("T2").Value = 0.009586633444024011

# Row 3
$ws.Range("O3").Value = 0.8587907398420774
$ws.Range("P3").Value = 0.8587907398420773
$ws.Range("S3").Value = 0.8365474849530689
$ws.Range("T3").Value = 0.8365474849530689

# Row 4
$ws.Range("O4").Value = 0.1313677243502452
$ws.Range("P4").Value = 0.1313677243502452
$ws.Range("S4").Value = 0.1279652123745701
$ws.Range("T4").Value = 0.1279652123745701

# Row 5
$ws.Range("M5").Value = 3.456265333333333
$ws.Range("N5").Value = 10.368796
$ws.Range("O5").Value = 0.009841535807677501
$ws.Range("P5").Value = 0.0098415358076775
$ws.Range("Q5").Value = 0.02161202712933333
$ws.Range("R5").Value = 0.194508244164
$ws.Range("S5").Value = 0.00025490236365349
$ws.Range("T5").Value = 0.00025490236365349

# Row 6
$ws.Range("O6").Value = 0.8587907398420774
$ws.Range("P6").Value = 0.8587907398420773
$ws.Range("S6").Value = 0.02224325488900852
$ws.Range("T6").Value = 0.02224325488900851

# Row 7
$ws.Range("O7").Value = 0.1313677243502452
$ws.Range("P7").Value = 0.1313677243502452
$ws.Range("S7").Value = 0.003402511975675061
$ws.Range("T7").Value = 0.003402511975675061
